$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$r = $p1.Range

# Replace the paragraph's visible text (excluding the trailing paragraph
# mark) with the updated placeholder ID. This also folds away the extra
# trailing-space run, leaving a single run behind.
$textRange = $d.Range($r.Start, $r.End - 1)
$textRange.Text = "**ID__AFFARS_AF_PGI_5301__ID**"

# Update the paragraph formatting: new left indent + a (style-less) paragraph
# border on all four sides with 5-twip spacing from the text.
$pf = $p1.Range.ParagraphFormat
$pf.LeftIndent = 11.25

$bs = $pf.Borders
$bs.DistanceFromTop = 5
$bs.DistanceFromLeft = 5
$bs.DistanceFromBottom = 5
$bs.DistanceFromRight = 5
